# Daily Update Feb 08 2020
# Appends the new "Cruise Ship" / "Others" data point (Diamond Princess
# cruise ship cases, reported near Yokohama, Japan) as row 84 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New row 84: State/Region, Country, lat, long
$ws.Range("A84").Value = "Cruise Ship"
$ws.Range("B84").Value = "Others"
$ws.Range("C84").Value = 35.4437
$ws.Range("D84").Value = 139.638

# Leave the cursor where the author's saved view shows it: one row below the
# newly appended data, in column C.
$ws.Range("C85").Select()

# Best-effort: scroll the view down to roughly match the author's saved
# topLeftCell ("A73"); harmless no-op if the runtime does not persist window
# scroll position on save.
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
